# subindo funcionalidade de delete
# Cronograma.xlsx - reorder/replace upcoming deliverables on the schedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: now shows the "search" deliverable (previously in row 10)
$ws.Range("B9").Value = "Implementação de busca para alteração de dados das classes"

# Row 10: now shows the "update" deliverable (previously in row 9)
$ws.Range("B10").Value = "Implementação das funcionalidades de atualização dos dados das classes"

# Row 11: the old "testes de homologação relacionados a buscas..." deliverable is
# removed entirely; the "delete" deliverable (previously row 13) moves up to row 11.
$ws.Range("B11").Value = "Implementação da funcionalidade de deletar um objeto de uma classe"

# Row 12 keeps its original text (index shifts internally, content unchanged).
$ws.Range("B12").Value = "Testes de homologação relacionados ao software em geral, para identificação de bugs e erros e correção dos mesmos"

# Row 13: brand new deliverable about renewing membership cards / promoting a student.
$newline = [char]10
$leftQuote = [char]8220
$rightQuote = [char]8221
$b13Text = "Mostrar a renovação da carteira de filiação de um aluno/professor com data vencida" + $newline + $leftQuote + "Promoção" + $rightQuote + " de um aluno a professor"
$ws.Range("B13").Value = $b13Text
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4107
$ws.Rows.Item(13).RowHeight = 45

# Row 14 keeps its original text (index shifts internally, content unchanged).
$ws.Range("B14").Value = "Homologação final do software, visando a correção de todos os erros e bugs identificados"

# Selection cursor moved to F12 in the saved file.
$ws.Range("F12").Select()
